# Applies the "zipLocation" UML business-rule / class-description / association /
# denormalization edits described by the commit:
#   "Updated UML to include zipLocation so that we don't have to roll down Person
#    attributes into each categorization. Also got rid of useless documents and
#    updated project specs doc to include class description and association
#    descriptions of zipLocation"

$d = $word.ActiveDocument

# Find the (single) paragraph whose visible text equals $text exactly (ignoring the
# trailing paragraph-mark newline). Anchoring by text, rather than a hard-coded
# Paragraphs(N) index, keeps the script robust to any earlier edits shifting indices.
function Find-ParaByText($doc, [string]$text) {
    foreach ($p in $doc.Paragraphs) {
        $t = $p.Range.Text
        $t = $t.TrimEnd([char]13, [char]10)
        if ($t -eq $text) {
            return $p
        }
    }
    return $null
}

# Replace a whole paragraph's contents (pPr + runs) with a literal WordprocessingML
# fragment - this is how list numbering (pStyle/numPr), proofErr spell-check markers,
# and multi-run splits get reproduced faithfully.
function Set-ParaXml($para, [string]$xmlFrag) {
    $xml = @"
$xmlFrag
"@
    $null = $para.Range.InsertXML($xml)
}

# --- 1. "Denormalization" section rewrite -------------------------------------------
# 1a. Replace the old generic "4 groups of people" paragraph with the new
#     zipcode/subkey/zipLocation/varchar denormalization description.
$pCurr = Find-ParaByText $d "Currently our class diagram has the 4 groups of people the hospital is most dependent on as specializations of a generic Person class. In the Phase 1 final draft, we will remove this Person class and have all of it’s attributes in each group, instead of them migrating from Person to each specialization as foreign key attributes."
Set-ParaXml $pCurr @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/><w:t xml:space="preserve">Currently our class diagram has the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>zipcode</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, city, and state information factored out of the Person class because it is a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>subkey</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, where the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>zipcode</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> functionally determines the city and state. To </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>den</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>ormalize</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> our class diagram, we will merge the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>zipLocation</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> class into the Person class and have the address as one big </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>varchar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p>
'@

# 1b. Wrap the "Denormalization" heading run in proofErr spellStart/spellEnd.
$pDenorm = Find-ParaByText $d "Denormalization"
Set-ParaXml $pDenorm @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>Denormalization</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@

# --- 2. New Associations bullet: "A person lives in one and only one zipLocation..." -
$pHeadTechAssoc = Find-ParaByText $d "A technician in a lab is the head technician of none or more laboratories. A laboratory has one and only one head technician."
$null = $pHeadTechAssoc.Range.InsertParagraphAfter()
$pPerson = $pHeadTechAssoc.Next()
Set-ParaXml $pPerson @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">A person lives in one and only one </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>zipLocation</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. A </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>zipLocation</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> contains zero or many Persons.</w:t></w:r></w:p>
'@

# --- 3. Move <w:lastRenderedPageBreak/> up one bullet (visit -> surgeon bullet loses it)
$pVisit = Find-ParaByText $d "A visit results in none or more prescriptions. A prescription is composed of a visit."
Set-ParaXml $pVisit @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>A visit results in none or more prescriptions. A prescription is composed of a visit.</w:t></w:r></w:p>
'@
$pSurgeonAssoc = Find-ParaByText $d "A surgeon performs none or more surgeries. A surgery is composed of a surgeon."
Set-ParaXml $pSurgeonAssoc @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>A surgeon performs none or more surgeries. A surgery is composed of a surgeon.</w:t></w:r></w:p>
'@

# --- 4. New class-description bullet: "ZipLocation - An integer code..." -------------
$pHeadTech = Find-ParaByText $d "Head technician – A technician who is placed in charge of one or more laboratories that he/she works at."
$null = $pHeadTech.Range.InsertParagraphAfter()
$pZip = $pHeadTech.Next()
Set-ParaXml $pZip @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>ZipLocation</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> – An integer code that determines a city within a state.</w:t></w:r></w:p>
'@

# --- 5. Collapse the "Surgery" bullet's two runs (split by a stray _GoBack bookmark) -
$pSurgery = Find-ParaByText $d "Surgery – A procedure performed on patients by surgeons."
Set-ParaXml $pSurgery @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Surgery – A procedure performed on patients by surgeons.</w:t></w:r></w:p>
'@

Write-Host "Done. Paragraphs: $($d.Paragraphs.Count)"
